$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.568.37'
$ws.Range("E2").Value = '  +0.29%  '
# Row 3
$ws.Range("D3").Value = '1.923.47'
$ws.Range("E3").Value = '  +0.50%  '
# Row 4
$ws.Range("E4").Value = '  +0.48%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.27'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.17%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.012'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.48%  '
# Row 7
$ws.Range("E7").Value = '  -0.10%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4054'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.19%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08219'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.90%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.011'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.10%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.84'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.36%  '
# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.112'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.77%  '
# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.904.04'
$ws.Range("E13").Value = '  +0.44%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.295'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.26%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.77'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.80%  '
# Row 16
$ws.Range("E16").Value = '  +1.49%  '
# Row 17
$ws.Range("E17").Value = '  +0.48%  '
# Row 18
$ws.Range("E18").Value = '  +0.18%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.62'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.36%  '
# Row 20
$ws.Range("E20").Value = '  +0.42%  '
# Row 21
$ws.Range("D21").Value = '29.572.45'
$ws.Range("E21").Value = '  +0.28%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.676'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.94%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.99'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.93%  '
# Row 25
$ws.Range("D25").Value = '2.157.06'
$ws.Range("E25").Value = '  +1.14%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.19'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.28%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.399'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.44%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.05'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.16%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.092'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.85%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.74'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.77%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.013'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.05%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09610'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.85%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.608'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.46%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.569'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.04%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.379'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.61%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06366'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.43%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02293'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.16%  '
# Row 38
$ws.Range("E38").Value = '  +1.64%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5955'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.24%  '
# Row 40
$ws.Range("E40").Value = '  +0.39%  '
# Row 41
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.012'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.46%  '
# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.867'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.46%  '
# Row 43
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1851'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.13%  '
# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.433'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.36%  '
# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.246'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.92%  '
# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.38'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.33%  '
# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07548'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.05%  '
# Row 48
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5559'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.27%  '
# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.994'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.85%  '
# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '119.55'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.14%  '
# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.440'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.93%  '
